$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values for the rows that changed
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "275.37"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.92"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.366"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06252"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.668"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.682"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.365"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8305"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01379"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1636"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08333"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03421"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03086"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09314"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.880"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001638"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04769"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006331"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005567"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001089"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.718"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.370"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3377"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04714"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007037"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003456"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01200"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006266"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8993"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03244"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002298"

# Update Hora (column G) from 8 to 9 for every data row (2-51)
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "9"
